# Applies the "fixed workflow" re-run: the first 4 cutoff rows (Cutoff = 0..3)
# are dropped from each results sheet, and the remaining rows shift up so the
# table now spans A1:C16 instead of A1:C20. The Cutoff column (A) is
# re-numbered sequentially (0..14) while the Reaction_number data (B/C)
# simply travels up with its original row.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Remove the first 4 data rows (rows 2-5); Excel shifts rows 6-20 up to 2-16.
    $ws.Range("A2:A5").EntireRow.Delete() | Out-Null

    # Re-sequence the Cutoff column (A) starting back at 0 for the new row 2.
    for ($i = 0; $i -le 14; $i++) {
        $ws.Cells.Item($i + 2, 1).Value = $i
    }
}
